# Update crypto price/volume figures per the Jan 6 2023 GitHub Actions refresh commit.
# Source cells are plain text (t="inlineStr") holding numbers/percentages formatted as
# strings (e.g. "259.18", "0.90%"); set NumberFormat to Text ("@") before writing so
# Excel stores the literal text instead of auto-coercing to a number/percentage.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '259.18'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '0.90%'
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '26.95'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-0.12%'
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '4.683'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '1.38%'
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.06047'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '2.83%'
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.660'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '0.64%'
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.8601'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '0.15%'
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9215'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-2.07%'
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1396'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-0.92%'
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.05171'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '24.27%'
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07089'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-0.20%'
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03060'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-3.86%'
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09133'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.22%'
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001530'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-1.14%'
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0006073'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '0.54%'
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.006070'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-2.39%'
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.470'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-1.62%'
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.168'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-1.19%'
# Row 19
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-0.95%'
# Row 20
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '2.45%'
# Row 21
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-0.19%'
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.113'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '7.57%'
# Row 23
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '0.10%'
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001217'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-0.40%'
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004019'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-6.41%'
# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-0.05%'
# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '-21.35%'
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03855'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '0.63%'
# Row 41
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '1.13%'
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.004038'
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.01497'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '30.71%'
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002199'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-9.55%'
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005195'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-5.06%'
# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-0.05%'
# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '6.47%'
# Row 48
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-38.55%'
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002099'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-0.05%'
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0001999'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.05%'
